# First test for User Part — add tester feedback ("programmer check" column E,
# "tester feedback" column G) for the register / login / facility / booking
# rows, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen the feedback column (G) so the new, longer comments are readable.
$ws.Columns.Item(7).ColumnWidth = 55.83203125

# --- Registration (row 13) ---------------------------------------------
$ws.Range("E13").Value = "check"

$fb13 = @'
S2
1. email address should be unique, but I used same email to register, register function could work. That's not right.
2. when I input an inlegal email address, the box border will be red(that's right), but I suggest that it should show some reminder(shown on the page or some alert box) to remind user to input a valid email. That will be better.
3. I have not recieved any confirmation code when I register(I know it has not completed, I just mention it)
4.the password should be inputted twice for double check.
5.the password should be 6-20 length.
6.the password should only contain letters and numbers, but when I input  this:"," it also could work.
7.just suggestion on the page:
7.1 "ID", this word, may make confuse to users, it can be changed to "username"
7.2 "name" may be divided into "first name" and "last name". That will be better I think.
'@
$ws.Range("G13").Value = $fb13
$ws.Range("G13").WrapText = $true
$ws.Range("G13").HorizontalAlignment = -4131

# --- Login (row 14) ------------------------------------------------------
$ws.Range("E14").Value = "check"

$fb14 = @'
well done
I think when the register part is completed, the login part is fine.
'@
$ws.Range("G14").Value = $fb14
$ws.Range("G14").WrapText = $true

# --- Account recovery / Account info update (rows 15-16) -----------------
$notDone = "it has not completed, right?"
$ws.Range("E15").Value = $notDone
$ws.Range("E16").Value = $notDone

# --- Search facility (row 18) --------------------------------------------
$ws.Range("E18").Value = "check"

$fb18 = @'

well done

'@
$ws.Range("G18").Value = $fb18
$ws.Range("G18").WrapText = $true

# --- Book facility/event (row 23) - mixed formatting feedback -----------
$ws.Range("E23").Value = "check"

$fb23Part1 = @'
S2
1. no confirmation email
2. how can I prove that I am a member of the University? How can I get the discount?
3. when I booked a facility, the time is not right:the start time is 4pm and the end time is 9am.And I can book the day before today, like 9th May. It is not logic.
4. I could not see what exact time I have booked(there is just the end time).
5. (confuse) I just wonder that is it right that I can see other members booking on this calendar as an user.
6. as an user, I can cancel other member's booking? That is not right.

'@
$fb23Part2 = "7. the 3) 4) I have not tested, I want to test it after the booking part is finished."
$fb23 = $fb23Part1 + $fb23Part2 + "`n"

$ws.Range("G23").Value = $fb23
$ws.Range("G23").WrapText = $true

$redStart = $fb23Part1.Length + 1
$redLen = $fb23Part2.Length
$redChars = $ws.Range("G23").Characters($redStart, $redLen)
$redChars.Font.Size = 12
$redChars.Font.Color = 255
$redChars.Font.Name = "Calibri (Body)"

# --- Add (block) booking (row 27) ----------------------------------------
$ws.Range("E27").Value = "check"

$fb27 = @'
well done.
I will test it again after the booking part is finished.
'@
$ws.Range("G27").Value = $fb27
$ws.Range("G27").WrapText = $true

# Leave the selection where the author left it when finishing the review.
$ws.Range("G27").Select()
